# Updates cryptos list figures (price + 1h volume change) as scraped by the
# periodic GitHub Actions job, reflecting the latest pull from coinranking.com.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prices are stored as plain text in this sheet (not numbers), so force text
# entry (leading apostrophe) and reset the style afterwards so Excel doesn't
# silently reinterpret numeric-looking strings as numbers or leave a
# quote-prefix style behind.
function Set-TextValue($a1, $value) {
    $r = $ws.Range($a1)
    $r.Value = "'" + $value
    $r.Style = "Normal"
}

function Set-Row($row, $price, $volume) {
    if ($null -ne $price) {
        Set-TextValue "D$row" $price
    }
    if ($null -ne $volume) {
        $ws.Range("E$row").Value = $volume
    }
}

Set-Row 2  "64.333.15"   "  -2.86%  "
Set-Row 3  "3.168.00"    "  -4.50%  "
Set-Row 4  $null         "  -0.01%  "
Set-Row 5  "569.02"      "  -2.94%  "
Set-Row 6  "168.22"      "  -8.15%  "
Set-Row 7  "0.609"       "  -5.84%  "
Set-Row 9  "3.170.86"    "  -4.33%  "
Set-Row 10 $null         "  -4.56%  "
Set-Row 11 $null         "  -0.28%  "
Set-Row 12 "0.385"       "  -4.08%  "
Set-Row 13 "3.720.46"    "  -4.57%  "
Set-Row 14 $null         "  -2.41%  "
Set-Row 15 "64.375.22"   "  -2.90%  "
Set-Row 16 "25.33"       "  -3.75%  "
Set-Row 17 $null         "  -2.80%  "
Set-Row 18 "3.169.88"    "  -3.03%  "
Set-Row 19 "418.93"      "  -2.44%  "
Set-Row 20 "5.36"        "  -2.94%  "
Set-Row 21 "12.82"       "  -3.67%  "
Set-Row 22 "7.05"        "  -4.91%  "
Set-Row 23 "0.999"       "  -0.17%  "
Set-Row 24 "69.73"       "  -3.21%  "
Set-Row 25 $null         "  +0.43%  "
Set-Row 26 "0.484"       "  -6.38%  "
Set-Row 27 "0.0000105"   "  -7.32%  "
Set-Row 28 "8.85"        "  -1.30%  "
Set-Row 29 $null         "  +0.83%  "
Set-Row 30 $null         "  -6.01%  "
Set-Row 31 "21.64"       "  -3.28%  "
Set-Row 33 "5.01"        "  -3.60%  "
Set-Row 34 "6.31"        "  -4.45%  "
Set-Row 35 "1.13"        "  -4.61%  "
Set-Row 36 "157.23"      "  -1.67%  "
Set-Row 37 $null         "  -6.46%  "
Set-Row 38 "2.721.21"    "  -6.20%  "
Set-Row 39 $null         "  -6.64%  "
Set-Row 40 "24.25"       "  -8.91%  "
Set-Row 41 $null         "  -4.00%  "
Set-Row 42 $null         "  -2.52%  "
Set-Row 43 "0.708"       "  -7.51%  "
Set-Row 44 "0.0619"      "  -7.13%  "
Set-Row 45 "5.58"        "  -7.19%  "
Set-Row 46 $null         "  -3.86%  "
Set-Row 47 "292.99"      "  -7.28%  "
Set-Row 48 "21.53"       "  -7.54%  "
Set-Row 51 "0.0985"      $null

# Rows 49 and 50 swapped coins (FirstDigitalUSD now ranks above dogwifhat)
$ws.Range("B49").Value = "FirstDigitalUSD"
$ws.Range("C49").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D49" "1.00"
$ws.Range("E49").Value = "  -0.01%  "

$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D50" "2.01"
$ws.Range("E50").Value = "  -13.24%  "
